# Add a new worksheet "FI2CaCOCaCOSil" after "SalSIL"
$wb = $excel.ActiveWorkbook

$salSil = $wb.Worksheets.Item("SalSIL")
$newSheet = $wb.Worksheets.Add($null, $salSil)
$newSheet.Name = "FI2CaCOCaCOSil"

# Header row
$newSheet.Range("A1").Value = "FI2CaCOCaCOSil"

# Column labels
$newSheet.Range("A2").Value = "No."
$newSheet.Range("B2").Value = "M"
$newSheet.Range("C2").Value = "N"

# Data rows
$newSheet.Range("A3").Value = 317
$newSheet.Range("B3").Value = 0.58297928836962287
$newSheet.Range("C3").Value = 0.92219397238449263

$newSheet.Range("A4").Value = 318
$newSheet.Range("B4").Value = 0.58411519416210578
$newSheet.Range("C4").Value = 0.91104052645295797

$newSheet.Range("A5").Value = 319
$newSheet.Range("B5").Value = 0.57765977621847142
$newSheet.Range("C5").Value = 0.88765029395031281

# Make the new sheet the active/selected one
$newSheet.Select()
$newSheet.Range("A1:C2").Select()
